$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "edit2"
$ws.Range("B20").Value = "riya-morankar"
$ws.Range("C20").Value = "Merged"
$ws.Range("D20").Value = "N/A"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2025-06-19"
$ws.Range("F20").Value = "c886f687e6bc4f27615b5182d10b73894e43a993"
